# Apply the DTR report edits described by the commit diff:
#  - Append " R" to the shared remark text used by P15/P16
#    ("~OB Others|SIT Inbound|" -> "~OB Others|SIT Inbound| R")
#  - Move the half-day marker from "NO. OF HOURS LATE" (E) to
#    "NO. OF HOURS UNDERTIME" (F) for the two affected rows, and
#    correct the value in E9 from 0.75 to 0.25.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remarks column (P) - both rows referencing the same remark text.
$ws.Range("P15").Value = "~OB Others|SIT Inbound| R"
$ws.Range("P16").Value = "~OB Others|SIT Inbound| R"

# Row 7/8 block: clear E7, set F8.
$ws.Range("E7").Value = ""
$ws.Range("F8").Value = 0.25

# Row 9: 0.75 -> 0.25
$ws.Range("E9").Value = 0.25

# Row 17/18 block: clear E17, set F18.
$ws.Range("E17").Value = ""
$ws.Range("F18").Value = 0.25
